$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = 0
    9  = 2
    12 = -9
    15 = -5
    16 = -3
    17 = -6
    18 = -6
    20 = 0
    21 = -2
    22 = 4
    24 = 0
    27 = 2
    28 = -7
    29 = -4
    31 = -2
    34 = 3
    35 = -3
    36 = -5
    42 = -5
    45 = -5
    46 = 10
    48 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
